# Updated cryptos list on Mon Dec 25 10:50:25 UTC 2023 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns on Sheet1 with new
# crypto-ranking data. Price cells keep their original "Text" storage
# (the source data uses locale-formatted strings like "43.325.01" that
# are not valid numbers), so NumberFormat is forced to "@" (Text) before
# assigning each Price cell's value to avoid Excel auto-coercing them to
# numbers (which would also silently drop trailing zeros, e.g. 264.20 ->
# 264.2). Volume cells are plain percentage strings with padding spaces
# and are assigned as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.183.45'
$ws.Range("E2").Value = '  -1.11%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.276.65'
$ws.Range("E3").Value = '  -0.41%  '

$ws.Range("E4").Value = '  -0.24%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '111.88'
$ws.Range("E5").Value = '  +1.48%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '264.20'
$ws.Range("E6").Value = '  -1.02%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.642'
$ws.Range("E7").Value = '  +2.99%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.609'
$ws.Range("E9").Value = '  -0.71%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '46.65'
$ws.Range("E10").Value = '  -1.48%  '

$ws.Range("E11").Value = '  -0.67%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.26'
$ws.Range("E12").Value = '  +5.06%  '

$ws.Range("E14").Value = '  -1.96%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.619.24'
$ws.Range("E15").Value = '  -0.28%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.861'
$ws.Range("E16").Value = '  +2.28%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.274.04'
$ws.Range("E17").Value = '  -0.25%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.171.07'
$ws.Range("E18").Value = '  -0.79%  '

$ws.Range("E19").Value = '  -1.20%  '

$ws.Range("E20").Value = '  +1.43%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.19'
$ws.Range("E21").Value = '  +0.05%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.43'
$ws.Range("E22").Value = '  -0.72%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '234.34'
$ws.Range("E23").Value = '  +1.11%  '

$ws.Range("E24").Value = '  +3.79%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.35'
$ws.Range("E25").Value = '  -2.79%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.32'
$ws.Range("E27").Value = '  -2.06%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '41.23'
$ws.Range("E28").Value = '  -1.43%  '

$ws.Range("E29").Value = '  -1.35%  '

$ws.Range("E30").Value = '  -0.57%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '173.52'
$ws.Range("E31").Value = '  -1.37%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.48'
$ws.Range("E32").Value = '  +0.00%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0897'
$ws.Range("E33").Value = '  -2.75%  '

$ws.Range("E34").Value = '  +1.19%  '

$ws.Range("E35").Value = '  +3.58%  '

$ws.Range("E36").Value = '  +6.17%  '

$ws.Range("E37").Value = '  -0.22%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.87'
$ws.Range("E38").Value = '  +2.99%  '

$ws.Range("E39").Value = '  -3.14%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.58'
$ws.Range("E40").Value = '  +7.63%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '14.21'
$ws.Range("E41").Value = '  +4.62%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '74.94'
$ws.Range("E42").Value = '  +4.63%  '

$ws.Range("E43").Value = '  -2.45%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.09'
$ws.Range("E44").Value = '  -2.58%  '

$ws.Range("E45").Value = '  +0.04%  '

$ws.Range("E46").Value = '  -1.56%  '

$ws.Range("E47").Value = '  +4.85%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.53'
$ws.Range("E48").Value = '  -2.87%  '

$ws.Range("E49").Value = '  -0.88%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '100.40'
$ws.Range("E50").Value = '  -1.17%  '

$ws.Range("E51").Value = '  -2.35%  '
